# Update the instrument characteristic value on Sheet1
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("H9").Value = "5.6 k Ohms with 56 Ohm dampener"

# Add a new empty Sheet2 positioned right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Keep Sheet1 as the active/selected sheet (matches the source workbook)
$ws1.Activate()
